$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's B/C/D columns (email / loai / trangThai) are being reordered so
# that "email" (formerly column B) moves to column D, and "loai" / "trangThai"
# (formerly C / D) shift left into B / C. Reproduce that with a real
# column cut + insert, which carries cell values, cell styles AND the
# per-column width metadata along with the move (matching the diff, where
# the 20.125 width follows the email column to D and the 15.625 width
# follows the loai column to B, while the now-unused column C reverts to
# the sheet's default width).
$ws.Columns("B:B").Cut()
$ws.Columns("E:E").Insert()

# Column C (now holding the old "trangThai" column, which never had a
# custom width) should go back to the sheet default width instead of
# keeping a leftover explicit <col> entry from the shift above.
$ws.Columns("C:C").ClearFormats()

# The two mailto hyperlinks used to live on column B (test_email@gmail.com /
# admin@gmail.com); they now belong on column D where that text ended up.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:test_email@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:admin@gmail.com")

# Hyperlinks.Add reformats the cell with its own style; restore the plain
# "Hyperlink" cell style (same one the cells already used before the move)
# so D2/D3 end up styled exactly like B2/B3 were.
$ws.Range("D2").Style = "Hyperlink"
$ws.Range("D3").Style = "Hyperlink"

# Update the sheet's saved selection from E13 to G11.
$ws.Range("G11").Select()
